$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "38.791.44"
$ws.Range("E2").Value = "  +2.74%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.093.27"
$ws.Range("E3").Value = "  +2.62%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.09%  "

# Row 5 - BNB
$ws.Range("D5").Value = "228.28"
$ws.Range("E5").Value = "  +0.42%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +1.07%  "

# Row 7 - Solana
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.70"
$ws.Range("E7").Value = "  +2.01%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.02%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +2.12%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "0.0837"
$ws.Range("E10").Value = "  -0.04%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -0.91%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "2.403.46"
$ws.Range("E12").Value = "  +2.59%  "

# Row 13 - Chainlink
$ws.Range("E13").Value = "  +3.85%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "21.93"
$ws.Range("E14").Value = "  +4.33%  "

# Row 15 - Polygon
$ws.Range("D15").Value = "0.798"
$ws.Range("E15").Value = "  +3.67%  "

# Row 16 - Polkadot
$ws.Range("E16").Value = "  -0.12%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.095.86"
$ws.Range("E17").Value = "  +2.43%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "38.734.69"
$ws.Range("E18").Value = "  +2.64%  "

# Row 19 - Litecoin
$ws.Range("D19").Value = "71.73"
$ws.Range("E19").Value = "  +3.27%  "

# Row 20 - Uniswap
$ws.Range("D20").Value = "6.03"
$ws.Range("E20").Value = "  +2.00%  "

# Row 21 - ShibaInu
$ws.Range("E21").Value = "  +1.62%  "

# Row 22 - BitcoinCash
$ws.Range("D22").Value = "226.56"
$ws.Range("E22").Value = "  +1.27%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  -0.37%  "

# Row 24 - Toncoin
$ws.Range("D24").Value = "2.39"
$ws.Range("E24").Value = "  -1.72%  "

# Row 26 - Monero
$ws.Range("D26").Value = "170.54"
$ws.Range("E26").Value = "  +0.62%  "

# Row 27 - Cosmos
$ws.Range("D27").Value = "9.44"
$ws.Range("E27").Value = "  +0.86%  "

# Row 28 - Kaspa
$ws.Range("E28").Value = "  +6.99%  "

# Row 29 - ImmutableX
$ws.Range("D29").Value = "1.44"
$ws.Range("E29").Value = "  +12.12%  "

# Row 30 - EthereumClassic
$ws.Range("D30").Value = "19.16"
$ws.Range("E30").Value = "  +1.98%  "

# Row 31 - Stellar
$ws.Range("E31").Value = "  +0.88%  "

# Row 32 - WEMIXToken
$ws.Range("E32").Value = "  +4.60%  "

# Row 33 - Filecoin
$ws.Range("E33").Value = "  +2.63%  "

# Row 35 - Hedera
$ws.Range("E35").Value = "  +1.98%  "

# Row 36 - THORChain
$ws.Range("D36").Value = "6.42"
$ws.Range("E36").Value = "  -2.31%  "

# Row 37 - LidoDAOToken
$ws.Range("D37").Value = "2.38"
$ws.Range("E37").Value = "  +2.02%  "

# Row 38 - RenderToken
$ws.Range("E38").Value = "  +2.27%  "

# Row 39 - BinanceUSD
$ws.Range("E39").Value = "  -0.03%  "

# Row 40 - InjectiveProtocol
$ws.Range("D40").Value = "18.21"
$ws.Range("E40").Value = "  +1.29%  "

# Row 41 - now Aave (was Maker)
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "101.20"
$ws.Range("E41").Value = "  +3.81%  "

# Row 42 - now Maker (was Aave)
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "1.538.80"
$ws.Range("E42").Value = "  +0.77%  "

# Row 43 - VeChain
$ws.Range("D43").Value = "0.0222"
$ws.Range("E43").Value = "  +3.45%  "

# Row 44 - now HuobiToken (was Cronos)
$ws.Range("B44").Value = "HuobiToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D44").Value = "2.82"
$ws.Range("E44").Value = "  -0.84%  "

# Row 45 - now Cronos (was HuobiToken)
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").Value = "0.0924"
$ws.Range("E45").Value = "  +1.89%  "

# Row 46 - FraxShare
$ws.Range("D46").Value = "7.67"
$ws.Range("E46").Value = "  +8.14%  "

# Row 47 - TrustWalletToken
$ws.Range("E47").Value = "  +0.91%  "

# Row 48 - FTXToken
$ws.Range("D48").Value = "4.11"
$ws.Range("E48").Value = "  -2.20%  "

# Row 49 - ARBITRUM
$ws.Range("E49").Value = "  +2.53%  "

# Row 50 - MXToken
$ws.Range("E50").Value = "  +1.20%  "

# Row 51 - RocketPoolETH
$ws.Range("D51").Value = "2.290.28"
$ws.Range("E51").Value = "  +2.64%  "
